$d = $word.ActiveDocument

$d.Content.Find.Execute("MEDIUMBLOB", $false, $false, $false, $false, $false,
                         $true, 1, $false, "LONGBLOB", 2)
